$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.931.47"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "2.360.33"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.690"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.13"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.27%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.631"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +26.06%  "
$ws.Range("E10").Value = "  +5.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.42"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.25"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +21.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.54"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +19.79%  "
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "2.708.88"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.91"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.920"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.43%  "
$ws.Range("D18").Value = "2.354.54"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "43.868.59"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("E20").Value = "  +4.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.43"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.70"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.26"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +12.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.78"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +15.52%  "
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.21"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.81"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("E33").Value = "  +5.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.30"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0752"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.35"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.82"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.45"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0277"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.96"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.201"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +17.70%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.102"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.26%  "
$ws.Range("E46").Value = "  +4.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.51"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +13.36%  "
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.66"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.24%  "
